# Adding the changes we made on may 9th
# Accelerometer capture grew: 5 samples were recorded slightly earlier
# (now the new top of the series) and 5 more were captured at the end,
# so the whole x/y/z table grows from 20 to 30 data rows (A2:C21 -> A2:C31).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full x,y,z table for the new A2:C31 range (old A2:C21 data now sits at A7:C26).
$rows = @(
    @(-2.143470001220703, 2.89409122467041, 0.02666953206062343),
    @(-2.163377571105957, 2.95642032623291, 0.1989836648106574),
    @(-2.130269622802735, 2.900108051300049, 0.1112725704908369),
    @(-2.321649217605591, 2.889909553527832, -0.01538913398981086),
    @(-2.246681690216064, 2.861854076385498, 0.0640860199928283),
    @(-2.385556173324585, 2.812312602996826, -0.1285117015242577),
    @(-2.266226387023926, 2.866149425506591, 0.03211449086666163),
    @(-2.475775194168091, 2.822203063964844, -0.136699762195349),
    @(-2.403434085845947, 2.584281539916992, -1.100419497489932),
    @(-2.276368379592896, 2.15183401107788, -1.987968981266023),
    @(-2.559272384643557, 1.362055969238278, -2.492264986038209),
    @(-3.275786542892459, 0.1545626640319781, -2.795625829696655),
    @(-4.113391685485841, -0.9623365402221682, -2.931767368316651),
    @(-4.743673324584959, -0.9649306297302217, -3.303337550163267),
    @(-3.626681327819824, 0.6835670471191406, -2.28171181678772),
    @(-2.572325563430785, 0.671577787399292, -3.271707224845887),
    @(-6.29947957992554, 0.9659337997436523, -1.046118927001951),
    @(1.245738077163745, 0.1863647580146739, 1.624357903003705),
    @(0.1075388908385779, 2.011311626434347, 0.675841426849348),
    @(2.805893659591699, 2.456681728363033, -0.3629340231418602),
    @(4.720451354980446, 1.925141620635989, -0.02809072732925327),
    @(-0.6559253215789739, 2.539989399909972, 0.1421575546264644),
    @(0.5849607467651348, 2.491342687606815, -0.09362407922744952),
    @(0.230958747863774, 3.317881345748901, -0.5695523142814627),
    @(1.351141929626465, 3.224694967269897, -0.3476141095161438),
    @(0.705779266357419, 3.61713025569916, 0.1973931401968027),
    @(0.6601259231567385, 3.235907554626462, 0.04604268670081965),
    @(0.9756811141967799, 3.334283685684206, 0.05510960519313853),
    @(0.6400001525878865, 3.297549438476561, 0.02418547868728584),
    @(0.471990585327149, 3.205180048942565, -0.01620917022228236)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
